# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -3

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 7

$ws.Range("F6").Value = -2

$ws.Range("F7").Value = -1

$ws.Range("F10").Value = -1

$ws.Range("F12").Value = -1

$ws.Range("F17").Value = -3

$ws.Range("F18").Value = -4

$ws.Range("F19").Value = 4

$ws.Range("F22").Value = -5
